# The author repositioned Word's "last edit" bookmark (_GoBack) from around
# the diagram picture to the spot where they last typed/edited text: right
# after "...requests that were given as threads." in the "design decisions"
# paragraph. Word only ever keeps a single _GoBack bookmark, so re-adding it
# at the new location also removes it from its old location automatically.
#
# Along the way, the run that held that sentence gets split at the two
# points the cursor passed through while editing (between "req" and
# "uests", and between "threads" and the following period), which is why
# the diff shows that single run broken into three runs.

$d = $word.ActiveDocument

# Locate the boundary right after "...each of the req" (end of the first
# fragment of the split run).
$rFind1 = $d.Content
$rFind1.Find.Execute("each of the req", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0) | Out-Null
$splitPos1 = $rFind1.End

# Locate the boundary right after "...uests that were given as threads"
# (i.e. right before the trailing period).
$rFind2 = $d.Content
$rFind2.Find.Execute("uests that were given as threads", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0) | Out-Null
$splitPos2 = $rFind2.End

# Locate the boundary right after the period that follows "threads"
# (this is where the _GoBack bookmark needs to end up).
$rFind3 = $d.Content
$rFind3.Find.Execute("given as threads.", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0) | Out-Null
$splitPos3 = $rFind3.End

# Split the run in two places by dropping a temporary bookmark at each
# boundary and immediately deleting it again -- the bookmark insertion
# forces a run break at that character offset, and the break survives the
# bookmark's removal.
$d.Bookmarks.Add("__TmpSplit1", $d.Range($splitPos1, $splitPos1))
$d.Bookmarks("__TmpSplit1").Delete()

$d.Bookmarks.Add("__TmpSplit2", $d.Range($splitPos2, $splitPos2))
$d.Bookmarks("__TmpSplit2").Delete()

# Re-add _GoBack collapsed right after the period. Because a document can
# only have one bookmark named "_GoBack", this both creates it here and
# removes the old one that wrapped the diagram's paragraph.
$d.Bookmarks.Add("_GoBack", $d.Range($splitPos3, $splitPos3))
